# Update Name of Algo
# Applies updated RandomForest imputation results to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = -11.36579999999999
$ws.Range("C3").Value = -11.9707
$ws.Range("C5").Value = -11.68509999999999
$ws.Range("E7").Value = 12.1757
$ws.Range("A9").Value = -20.73829999999998
$ws.Range("E9").Value = 12.80340000000001
$ws.Range("C11").Value = -13.97860000000001
$ws.Range("C12").Value = -14.02250000000001
$ws.Range("A13").Value = -21.96900000000001
$ws.Range("A16").Value = -19.92079999999998
$ws.Range("A18").Value = -21.6683
$ws.Range("A20").Value = -22.05480000000003
$ws.Range("C21").Value = -14.04320000000001
$ws.Range("E21").Value = 12.94929999999999
